$d = $word.ActiveDocument

# The body paragraph holds a single run whose text is split into six
# sentences by pairs of manual line breaks (<w:br/><w:br/>). Word exposes
# each manual break as a vertical-tab character (chr 11) in Range.Text, so
# splitting on that character recovers the six editable text segments.
$para = $d.Paragraphs(2)
$r = $para.Range
$baseStart = $r.Start
$fulltext = $r.Text
$parts = $fulltext.Split([char]11)

$newParts = @(
  "The Gross Regional Product (GRP) serves as a vital economic barometer, capturing the total economic output of a region. In the case of San Diego County, the GRP data spanning from 2019 to 2023 offers a revealing glimpse into the area's economic performance and growth trajectory across various industries.",
  "Over this five-year period, San Diego County's GRP has consistently trended upwards. Beginning at approximately `$244.28 billion in 2019, it surged to about `$308.71 billion by 2023, marking a notable growth of around 26.4%. This upward trajectory is mirrored in the GRP per capita, which rose from `$73,346.92 in 2019 to `$94,915.87 in 2023, indicating a boost in economic productivity per individual within the region.",
  "Several industries have been pivotal in driving this growth. The government sector emerged as the largest contributor, with its economic input climbing from `$45.19 billion in 2019 to `$52.92 billion in 2023. Manufacturing also played a crucial role, with its contribution increasing from `$28.26 billion to `$31.67 billion over the same period. The Professional, Scientific, and Technical Services sector experienced substantial growth, expanding from `$28.13 billion to `$37.04 billion. Additionally, the Health Care and Social Assistance sector saw its contribution rise from `$15.38 billion to `$20.21 billion.",
  "The year 2020, marked by the COVID-19 pandemic, presented challenges, particularly for sectors such as Accommodation and Food Services, and Arts, Entertainment, and Recreation, which experienced slight declines. Nevertheless, sectors like Finance and Insurance, along with Government, demonstrated resilience and even growth during this tumultuous period.",
  "When compared to state and national levels, San Diego County's economic performance stands out. By 2023, the county's GRP per capita reached `$94,915.87, surpassing the California state average of `$93,799.67 and significantly exceeding the national average of `$77,366.43. This comparison underscores the county's robust economic health relative to broader benchmarks.",
  "In conclusion, San Diego County has exhibited impressive economic growth over the past five years, with substantial contributions from the government, manufacturing, and professional services sectors. The rise in GRP per capita reflects improved productivity and economic well-being for its residents. Despite the challenges posed by the pandemic, the region's economy has shown remarkable resilience and adaptability, positioning it favorably for continued future growth."
)

# Compute the character offset (within $fulltext) at which each split part
# begins, so we can target each sentence with a precise Range.
$offset = 0
$starts = @()
foreach ($p in $parts) {
    $starts += $offset
    $offset += $p.Length + 1
}

# Replace each of the six text sentences (even-indexed parts; odd-indexed
# parts are the single-char gaps where the breaks were split out) using
# direct Range.Text assignment rather than Find.Execute, since Find-based
# replacement triggers AutoCorrect's smart-quote substitution and would
# turn straight apostrophes into curly ones. Work from the last sentence
# to the first so earlier offsets stay valid as text lengths change.
for ($i = 5; $i -ge 0; $i--) {
    $partIndex = $i * 2
    $s = $baseStart + $starts[$partIndex]
    $e = $s + $parts[$partIndex].Length
    $sub = $d.Range($s, $e)
    $sub.Text = $newParts[$i]
}

Write-Host "Replacement complete"
